$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: update the "last modified" date field result
#   2020-09-03 -> 2021-02-01
# ------------------------------------------------------------------
$d.Content.Find.Execute("2020-09-03", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "2021-02-01", 2) | Out-Null

# ------------------------------------------------------------------
# Change 2: "<g type="leaf"/>" -> "<g ref="#leaf"/>"
#   (the sample row in the <g>/leaf symbol table)
# The final run needs to be split into five runs:
#   "<g "  |  "ref"  |  ="  |  #  |  leaf"/>
# all sharing the exact same run formatting as the original run.
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute('type="leaf"')
if (-not $found) {
    throw "Could not find the <g type=`"leaf`"/> sample text"
}

$rstart = $rng.Start
# Plain textual substitution first (stays inside the single found run).
$rng.Text = 'ref="#leaf"'

# Segment boundaries (relative offsets) within the text just inserted,
# chosen so the final runs read: ref | =" | # | leaf"/>
# (the leading "<g " stays part of the run that precedes $rstart).
$segments = @(
    @(0, 3),    # ref
    @(3, 5),    # ="
    @(5, 6),    # #
    @(6, 13)    # leaf"/>
)

foreach ($seg in $segments) {
    $segStart = $rstart + $seg[0]
    $segEnd   = $rstart + $seg[1]
    $segRange = $d.Range($segStart, $segEnd)
    # Toggling a character property on/off forces this sub-range to be
    # materialised as its own run without actually changing any
    # visible formatting (Bold ends up right back where it started).
    $segRange.Bold = 1
    $segRange.Bold = 0
}
